$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-13 (columns B and D changes, per diff)
$ws.Range("D2").Value = 0.3406929969787598
$ws.Range("D3").Value = 0.2575700283050537
$ws.Range("D4").Value = 0.009903192520141602
$ws.Range("D5").Value = 0.0007691383361816406
$ws.Range("D6").Value = 0.07412600517272949
$ws.Range("D7").Value = 0.03745579719543457

$ws.Range("B8").Value = 4
$ws.Range("D8").Value = 0.009923219680786133

$ws.Range("B9").Value = 16
$ws.Range("D9").Value = 4.33816385269165

$ws.Range("B10").Value = 16
$ws.Range("D10").Value = 2.093926191329956

$ws.Range("B11").Value = 4
$ws.Range("D11").Value = 0.00992894172668457

$ws.Range("B12").Value = 16
$ws.Range("D12").Value = 64.24477481842041

$ws.Range("B13").Value = 16
$ws.Range("D13").Value = 33.1779158115387

# Add new rows 14-19
$ws.Range("A14").Value = "Initial CFF"
$ws.Range("B14").Value = 3
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 0.2419688701629639

$ws.Range("A15").Value = "Grow CFF"
$ws.Range("B15").Value = 9
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 0.7095961570739746

$ws.Range("A16").Value = "Direct CFF"
$ws.Range("B16").Value = 9
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 0.2056238651275635

$ws.Range("A17").Value = "Initial CFF"
$ws.Range("B17").Value = 3
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 0.003067970275878906

$ws.Range("A18").Value = "Grow CFF"
$ws.Range("B18").Value = 9
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3.822492122650146

$ws.Range("A19").Value = "Direct CFF"
$ws.Range("B19").Value = 9
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 1.816632032394409
